$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Weekly Quantity" ---
$ws1 = $wb.Worksheets.Item("Weekly Quantity")

# Row 9 gets the values that used to be in row 12 (date 45130.99999999999, qty 470)
$ws1.Range("A9").Value = 45130.99999999999
$ws1.Range("B9").Value = 470

# Rows 10-13 are no longer part of the data range; delete the entire rows
# (a plain Range.Delete shifts/clobbers unexpectedly in this runtime).
$ws1.Range("A10:B13").EntireRow.Delete()

# --- Sheet 2: "Monthly Trend" ---
$ws2 = $wb.Worksheets.Item("Monthly Trend")

# Row 6 gets the date that used to be in row 7 (45138.99999999999) and a new qty value 470
$ws2.Range("A6").Value = 45138.99999999999
$ws2.Range("B6").Value = 470

# Row 7 is no longer part of the data range; delete the entire row
# (a plain Range.Delete shifts/clobbers unexpectedly in this runtime).
$ws2.Range("A7:B7").EntireRow.Delete()
